# crosswalk.xlsx: rename the lone sheet to "LTER", add a new "EULI" sheet
# after it with clarifying notes about the "other" phyto bucket (see
# commit message "clarify other phytos in crosswalk/script"), and make
# EULI the active tab/selection.

$wb = $excel.ActiveWorkbook

# Rename the existing (only) sheet "Sheet2" -> "LTER"
$lter = $wb.Worksheets.Item(1)
$lter.Name = "LTER"

# Add a brand-new worksheet right after LTER, named "EULI"
$euli = $wb.Worksheets.Add($null, $lter)
$euli.Name = "EULI"

# Populate EULI with the new clarifying notes (one per row, column A)
$euli.Range("A1").Value = "crypto - crypto"
$euli.Range("A2").Value = "chloro - chloro"
$euli.Range("A3").Value = "cyano cyano match"
$euli.Range("A4").Value = "Ochryophya classes (Bac, Frag, Cosc) = diatoms"
$euli.Range("A5").Value = "other can be average of everything else in PLoS paper"

# Leave the selection on A7 of the newly-active EULI sheet
[void]$euli.Range("A7").Select()
